$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Add new row of data (row 4): A=TestCase, B=Name, C=MobileNumber
# (write Name/MobileNumber/TestCase in this order so new shared-string
# entries land in the same index order as the target workbook)
$ws.Range("B4").Value = "Amira"
$ws.Range("C4").Value = "1116332215&"
$ws.Range("A4").Value = "03-edit Mobile Nuber while login "

# Left-align the new MobileNumber cell (C4), matching the added cellXf with horizontal="left"
$ws.Range("C4").HorizontalAlignment = -4131  # xlHAlignLeft

# Column A width / bestFit adjustment (widen to fit the new longer TestCase text)
$ws.Columns.Item(1).ColumnWidth = 41.6

# Update the active cell selection on the sheet view
$ws.Range("A9").Select() | Out-Null

$wb.Save()
